$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values
# (e.g. "14.00", "252.60") are not silently coerced into numbers,
# matching the original inline-string storage used by this sheet.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '31.483.17'
$ws.Cells.Item(2, 5).Value = '  +3.75%  '
$ws.Cells.Item(3, 4).Value = '1.987.88'
$ws.Cells.Item(3, 5).Value = '  +5.99%  '
$ws.Cells.Item(4, 4).Value = '0.9996'
$ws.Cells.Item(4, 5).Value = '  -0.15%  '
$ws.Cells.Item(5, 4).Value = '0.8032'
$ws.Cells.Item(5, 5).Value = '  +71.47%  '
$ws.Cells.Item(6, 4).Value = '252.60'
$ws.Cells.Item(6, 5).Value = '  +3.53%  '
$ws.Cells.Item(7, 4).Value = '0.9977'
$ws.Cells.Item(7, 5).Value = '  -0.39%  '
$ws.Cells.Item(8, 4).Value = '0.3435'
$ws.Cells.Item(8, 5).Value = '  +19.70%  '
$ws.Cells.Item(9, 4).Value = '25.77'
$ws.Cells.Item(9, 5).Value = '  +17.32%  '
$ws.Cells.Item(10, 4).Value = '0.06947'
$ws.Cells.Item(10, 5).Value = '  +8.12%  '
$ws.Cells.Item(11, 4).Value = '0.8374'
$ws.Cells.Item(11, 5).Value = '  +16.15%  '
$ws.Cells.Item(12, 4).Value = '0.08123'
$ws.Cells.Item(12, 5).Value = '  +4.29%  '
$ws.Cells.Item(13, 4).Value = '102.14'
$ws.Cells.Item(13, 5).Value = '  +7.42%  '
$ws.Cells.Item(14, 4).Value = '1.983.68'
$ws.Cells.Item(14, 5).Value = '  +5.75%  '
$ws.Cells.Item(15, 4).Value = '5.493'
$ws.Cells.Item(15, 5).Value = '  +6.91%  '
$ws.Cells.Item(16, 4).Value = '274.99'
$ws.Cells.Item(16, 5).Value = '  -1.35%  '
$ws.Cells.Item(17, 4).Value = '31.475.06'
$ws.Cells.Item(17, 5).Value = '  +3.79%  '
$ws.Cells.Item(18, 4).Value = '14.00'
$ws.Cells.Item(18, 5).Value = '  +8.09%  '
$ws.Cells.Item(19, 4).Value = '0.000007873'
$ws.Cells.Item(19, 5).Value = '  +6.44%  '
$ws.Cells.Item(20, 4).Value = '2.240.59'
$ws.Cells.Item(20, 5).Value = '  +5.01%  '
$ws.Cells.Item(21, 4).Value = '5.691'
$ws.Cells.Item(21, 5).Value = '  +9.13%  '
$ws.Cells.Item(22, 4).Value = '0.9946'
$ws.Cells.Item(22, 5).Value = '  -0.76%  '
$ws.Cells.Item(23, 4).Value = '0.9944'
$ws.Cells.Item(23, 5).Value = '  -0.66%  '
$ws.Cells.Item(24, 4).Value = '6.868'
$ws.Cells.Item(24, 5).Value = '  +10.05%  '
$ws.Cells.Item(25, 4).Value = '0.1592'
$ws.Cells.Item(25, 5).Value = '  +66.26%  '
$ws.Cells.Item(26, 4).Value = '9.669'
$ws.Cells.Item(26, 5).Value = '  +7.26%  '
$ws.Cells.Item(27, 4).Value = '165.76'
$ws.Cells.Item(27, 5).Value = '  +1.34%  '
$ws.Cells.Item(28, 4).Value = '19.74'
$ws.Cells.Item(28, 5).Value = '  +5.80%  '
$ws.Cells.Item(29, 4).Value = '2.216'
$ws.Cells.Item(29, 5).Value = '  +17.93%  '
$ws.Cells.Item(30, 4).Value = '1.558'
$ws.Cells.Item(30, 5).Value = '  +6.26%  '
$ws.Cells.Item(31, 4).Value = '1.353'
$ws.Cells.Item(31, 5).Value = '  +1.66%  '
$ws.Cells.Item(32, 4).Value = '4.559'
$ws.Cells.Item(32, 5).Value = '  +8.37%  '
$ws.Cells.Item(33, 4).Value = '4.324'
$ws.Cells.Item(33, 5).Value = '  +5.75%  '
$ws.Cells.Item(34, 4).Value = '0.05203'
$ws.Cells.Item(34, 5).Value = '  +8.26%  '
$ws.Cells.Item(35, 4).Value = '1.218'
$ws.Cells.Item(35, 5).Value = '  +8.96%  '
$ws.Cells.Item(36, 4).Value = '0.7461'
$ws.Cells.Item(36, 5).Value = '  +8.92%  '
$ws.Cells.Item(37, 4).Value = '2.781'
$ws.Cells.Item(37, 5).Value = '  +2.64%  '
$ws.Cells.Item(38, 4).Value = '0.9927'
$ws.Cells.Item(38, 5).Value = '  -0.84%  '
$ws.Cells.Item(39, 4).Value = '0.01986'
$ws.Cells.Item(40, 4).Value = '2.911'
$ws.Cells.Item(40, 5).Value = '  +3.61%  '
$ws.Cells.Item(41, 4).Value = '6.595'
$ws.Cells.Item(41, 5).Value = '  +5.92%  '
$ws.Cells.Item(42, 4).Value = '78.47'
$ws.Cells.Item(42, 5).Value = '  +5.56%  '
$ws.Cells.Item(43, 4).Value = '0.4647'
$ws.Cells.Item(43, 5).Value = '  +10.20%  '
$ws.Cells.Item(44, 4).Value = '2.073'
$ws.Cells.Item(44, 5).Value = '  +7.37%  '
$ws.Cells.Item(45, 2).Value = 'Quant'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(45, 4).Value = '105.78'
$ws.Cells.Item(45, 5).Value = '  +5.06%  '
$ws.Cells.Item(46, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(46, 4).Value = '0.8534'
$ws.Cells.Item(46, 5).Value = '  +3.73%  '
$ws.Cells.Item(47, 4).Value = '0.9973'
$ws.Cells.Item(47, 5).Value = '  -0.37%  '
$ws.Cells.Item(48, 4).Value = '9.928'
$ws.Cells.Item(48, 5).Value = '  +3.70%  '
$ws.Cells.Item(49, 4).Value = '7.507'
$ws.Cells.Item(49, 5).Value = '  +8.56%  '
$ws.Cells.Item(50, 4).Value = '36.48'
$ws.Cells.Item(50, 5).Value = '  +4.26%  '
$ws.Cells.Item(51, 4).Value = '0.4257'
$ws.Cells.Item(51, 5).Value = '  +9.36%  '

# Restore the default (unstyled) cell style now that the text type is locked in,
# so the saved workbook does not pick up a stray number-format style change.
$priceRange.Style = "Normal"

